$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("advanced tasks")
$dst = $wb.Worksheets.Item("Sheet3")

# 1) Before touching "advanced tasks", duplicate its current A1:E7 table
#    (values + formatting) onto the still-empty "Sheet3" so the "due date"
#    column (and the strings/dates it used) keep living somewhere.
$src.Range("A1:E7").Copy()
$dst.Range("A1").PasteSpecial(-4163)
$src.Range("A1:E7").Copy()
$dst.Range("A1").PasteSpecial(-4122)

$dst.Columns.Item(3).ColumnWidth = 23.053385416666668
$dst.Columns.Item(4).ColumnWidth = 41.498697916666664
$dst.Columns.Item(5).ColumnWidth = 12.830729166666666

$excel.Goto($dst.Range("E4:E7"))

# 2) Remove the "due date" column (column B) from "advanced tasks" - this is
#    the actual edit described by the commit message.
$src.Range("B:B").Delete()

# Leave "advanced tasks" as the active sheet/selection, matching the final
# workbook view state.
$excel.Goto($src.Range("D4:D7"))
